$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting_id column (B) populated for both task rows.
$ws.Range("B2").Value = "MOM-20260107-001"
$ws.Range("B3").Value = "MOM-20260107-001"

# created_on (F) and deadline (K) switch from numeric/date-formatted values
# to plain text strings. Force text number format first so Excel does not
# re-parse the strings back into serial date numbers, then strip the
# leftover custom date format so the cells end up with the default
# (unformatted) style, matching the cleaned-up registry.
$ws.Range("F2:F3").NumberFormat = "@"
$ws.Range("F2").Value = "2026-01-07 19:51:37"
$ws.Range("F3").Value = "2026-01-07 19:51:37"

$ws.Range("K2:K3").NumberFormat = "@"
$ws.Range("K2").Value = "2026-01-14"
$ws.Range("K3").Value = "2026-01-14"

$ws.Range("F2:F3").ClearFormats()
$ws.Range("K2:K3").ClearFormats()
